$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# --- Row 2-44, 47-51: update Price (D) and Volume(1h) (E) columns ---
Set-TextValue "D2" "25.506.31"
Set-TextValue "E2" "  +2.28%  "
Set-TextValue "D3" "1.667.05"
Set-TextValue "E3" "  +1.86%  "
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  +0.39%  "
Set-TextValue "D5" "233.58"
Set-TextValue "E5" "  +0.40%  "
Set-TextValue "D6" "1.002"
Set-TextValue "E6" "  +0.25%  "
Set-TextValue "D7" "0.4598"
Set-TextValue "E7" "  -3.17%  "
Set-TextValue "D8" "0.2566"
Set-TextValue "E8" "  -0.30%  "
Set-TextValue "D9" "0.06108"
Set-TextValue "E9" "  +0.39%  "
Set-TextValue "D10" "1.667.74"
Set-TextValue "E10" "  +1.94%  "
Set-TextValue "D11" "0.06950"
Set-TextValue "E11" "  -0.74%  "
Set-TextValue "D12" "14.56"
Set-TextValue "E12" "  -0.15%  "
Set-TextValue "D13" "4.325"
Set-TextValue "E13" "  -0.48%  "
Set-TextValue "D14" "74.71"
Set-TextValue "E14" "  +1.66%  "
Set-TextValue "D15" "0.5618"
Set-TextValue "E15" "  -6.04%  "
Set-TextValue "D16" "1.002"
Set-TextValue "E16" "  +0.20%  "
Set-TextValue "D17" "1.002"
Set-TextValue "E17" "  +0.45%  "
Set-TextValue "D18" "25.509.58"
Set-TextValue "E18" "  +2.38%  "
Set-TextValue "D19" "0.000006669"
Set-TextValue "E19" "  +1.54%  "
Set-TextValue "D20" "11.31"
Set-TextValue "E20" "  +1.46%  "
Set-TextValue "D21" "1.880.60"
Set-TextValue "E21" "  +1.54%  "
Set-TextValue "D22" "4.410"
Set-TextValue "E22" "  +1.68%  "
Set-TextValue "D23" "8.677"
Set-TextValue "E23" "  +1.72%  "
Set-TextValue "D24" "5.189"
Set-TextValue "E24" "  -0.63%  "
Set-TextValue "D25" "136.10"
Set-TextValue "E25" "  +2.42%  "
Set-TextValue "D26" "14.85"
Set-TextValue "E26" "  +0.28%  "
Set-TextValue "D27" "1.376"
Set-TextValue "E27" "  -0.40%  "
Set-TextValue "D28" "104.02"
Set-TextValue "E28" "  +0.64%  "
Set-TextValue "D29" "1.693"
Set-TextValue "E29" "  +3.70%  "
Set-TextValue "D30" "3.955"
Set-TextValue "E30" "  +2.32%  "
Set-TextValue "D31" "0.07730"
Set-TextValue "E31" "  +0.52%  "
Set-TextValue "D32" "3.592"
Set-TextValue "E32" "  +1.47%  "
Set-TextValue "D33" "0.04266"
Set-TextValue "E33" "  -0.86%  "
Set-TextValue "E34" "  +1.49%  "
Set-TextValue "D35" "0.9414"
Set-TextValue "E35" "  +2.01%  "
Set-TextValue "D36" "0.5941"
Set-TextValue "E36" "  +2.64%  "
Set-TextValue "D37" "0.9215"
Set-TextValue "E37" "  +13.62%  "
Set-TextValue "D38" "2.510"
Set-TextValue "E38" "  -0.70%  "
Set-TextValue "D39" "1.002"
Set-TextValue "E39" "  +0.36%  "
Set-TextValue "D40" "101.17"
Set-TextValue "E40" "  +3.93%  "
Set-TextValue "D41" "0.01454"
Set-TextValue "E41" "  -4.78%  "
Set-TextValue "D42" "1.813"
Set-TextValue "E42" "  +3.03%  "
Set-TextValue "D43" "0.3677"
Set-TextValue "E43" "  -0.13%  "
Set-TextValue "D44" "4.898"
Set-TextValue "E44" "  +4.23%  "
Set-TextValue "D47" "6.097"
Set-TextValue "E47" "  +1.07%  "
Set-TextValue "D48" "29.54"
Set-TextValue "E48" "  +0.44%  "
Set-TextValue "D49" "7.367"
Set-TextValue "E49" "  +2.50%  "
Set-TextValue "D50" "1.003"
Set-TextValue "E50" "  +0.39%  "
Set-TextValue "D51" "0.9992"
Set-TextValue "E51" "  +0.48%  "

# --- Rows 45 and 46: Algorand/Cronos swapped places in the ranking ---
Set-TextValue "B45" "Cronos"
Set-TextValue "C45" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D45" "0.05253"
Set-TextValue "E45" "  +1.34%  "

Set-TextValue "B46" "Algorand"
Set-TextValue "C46" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D46" "0.1099"
Set-TextValue "E46" "  +1.00%  "
